$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updates to columns F (AgTests) and G (AgPosit) for rows 423-455
# as part of the "Updated: pi 04. 06. 2021" commit.

$ws.Range("F423").Value = 439056

$ws.Range("F424").Value = 266182

$ws.Range("F425").Value = 138486
$ws.Range("G425").Value = 550

$ws.Range("F426").Value = 107049

$ws.Range("G427").Value = 367

$ws.Range("F431").Value = 171273
$ws.Range("G431").Value = 402

$ws.Range("F432").Value = 123613
$ws.Range("G432").Value = 428

$ws.Range("F433").Value = 86898
$ws.Range("G433").Value = 269

$ws.Range("F434").Value = 78577
$ws.Range("G434").Value = 278

$ws.Range("F435").Value = 82933

$ws.Range("F436").Value = 144868

$ws.Range("F437").Value = 166991
$ws.Range("G437").Value = 272

$ws.Range("F442").Value = 70256

$ws.Range("F446").Value = 86182
$ws.Range("G446").Value = 262

$ws.Range("F447").Value = 67401

$ws.Range("F448").Value = 61505
$ws.Range("G448").Value = 139

$ws.Range("F449").Value = 59910

$ws.Range("F450").Value = 91325
$ws.Range("G450").Value = 169

$ws.Range("F451").Value = 85532

$ws.Range("F452").Value = 74582

$ws.Range("F453").Value = 70017

$ws.Range("F454").Value = 51557
$ws.Range("G454").Value = 130

$ws.Range("F455").Value = 50303
$ws.Range("G455").Value = 117
